# Updated cryptos list price/volume figures (GitHub Actions scheduled refresh).
# Column D ("Price") and Column E ("Volume(1h)") values are plain text cells
# (stored as inline strings in the source workbook), so numeric-looking
# prices are written with a leading apostrophe to keep Excel from coercing
# them into real numbers (which would lose trailing zeros / exact text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.935.51"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.817.81"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'310.14"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4646"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'0.3705"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "'0.07358"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'0.8720"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "'20.46"
$ws.Range("D12").Value = "1.821.79"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'5.355"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "'0.07097"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "'6.515"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "'91.62"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'0.000008724"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "26.954.01"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "2.065.11"
$ws.Range("D25").Value = "'1.901"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'152.09"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "'2.146"
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("D29").Value = "'5.315"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'115.61"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "'0.08909"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'0.7599"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'1.155"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").Value = "'4.466"
$ws.Range("D35").Value = "'2.923"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.01959"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "'0.05262"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'2.944"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").Value = "'7.255"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.5342"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "'2.374"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D45").Value = "'8.467"
$ws.Range("D46").Value = "'0.4947"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "'1.681"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'103.36"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").Value = "'0.06286"
$ws.Range("E51").Value = "  -0.73%  "
